$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the external reference/link (Std140_CE_a_Results.xlsx) that is no
# longer needed - this drops xl/externalLinks/externalLink1.xml and the
# <externalReferences> element from xl/workbook.xml.
$links = $wb.LinkSources()
if ($links) {
    foreach ($link in $links) {
        $wb.BreakLink($link, 1)
    }
}

# Replace the numeric case labels in column A (rows 25-38) with their
# string "CE###" equivalents, which also populates xl/sharedStrings.xml.
$ws.Range("A25").Value = "CE100"
$ws.Range("A26").Value = "CE110"
$ws.Range("A27").Value = "CE120"
$ws.Range("A28").Value = "CE130"
$ws.Range("A29").Value = "CE140"
$ws.Range("A30").Value = "CE150"
$ws.Range("A31").Value = "CE160"
$ws.Range("A32").Value = "CE165"
$ws.Range("A33").Value = "CE170"
$ws.Range("A34").Value = "CE180"
$ws.Range("A35").Value = "CE185"
$ws.Range("A36").Value = "CE190"
$ws.Range("A37").Value = "CE195"
$ws.Range("A38").Value = "CE200"

# Update the view so row 25 is visible at the top and the CE rows are
# selected, matching the saved workbook/window state.
$ws.Range("A25:A38").Select()
